# This script applies an updated data-refresh snapshot to the "想去人数"
# (interest count) column F on each worksheet of the 上海-漫展信息 workbook,
# matching the values produced by the scraper run at commit 456a3b4.
#
# NOTE: this engine's PowerShell interpreter does not bind named
# (-Param value) arguments to custom functions, so positional arguments
# are used throughout.

$wb = $excel.ActiveWorkbook

function Set-Cells($SheetName, $Updates) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Updates.Keys) {
        $ws.Range($cellRef).Value = $Updates[$cellRef]
    }
}

# 展览 (Exhibitions)
Set-Cells "展览" @{
    "F2"  = 1719
    "F3"  = 10040
    "F6"  = 588
    "F14" = 469
    "F15" = 1165
    "F19" = 80
    "F22" = 312
    "F24" = 824
    "F31" = 320
    "F35" = 535
    "F38" = 1264
    "F39" = 802
}

# 演出 (Performances)
Set-Cells "演出" @{
    "F6"  = 86
    "F12" = 23
    "F13" = 94
    "F16" = 130
    "F19" = 512
    "F21" = 316
    "F23" = 66
    "F25" = 21
    "F30" = 203
    "F35" = 41
    "F36" = 29
    "F38" = 127
}

# 本地生活 (Local Life)
Set-Cells "本地生活" @{
    "F4"  = 800
    "F6"  = 2493
    "F7"  = 4011
    "F10" = 260
    "F11" = 169
}

# 全部类型 (All Types)
Set-Cells "全部类型" @{
    "F2"  = 1719
    "F4"  = 800
    "F5"  = 10040
    "F8"  = 4011
    "F10" = 260
    "F11" = 260
    "F12" = 588
    "F17" = 469
    "F18" = 1165
    "F27" = 312
    "F31" = 66
    "F39" = 535
    "F40" = 203
    "F43" = 802
    "F46" = 41
}
